$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9607400631401491, 0.6363141226038453, 1.2599803579699587, 2.757477484039489)"
$ws.Range("C2").Value = "NIG(1.3893131935033092, 1.0811932543645453, 3.171941863539942, 5.247304399970877)"
$ws.Range("D2").Value = "NCT(2.146763212066107, 1.4392517279164858, -0.6671619270341693, 2.229515124041949)"
$ws.Range("E2").Value = "JSU(-1.3597500685442423, 1.1992472612259442, 2.3616481091997383, 4.5831262051533965)"
